$wb = $excel.ActiveWorkbook

# --- Test Sheet 1: add a new row (row 10) with numeric data ---
$ws1 = $wb.Worksheets.Item("Test Sheet 1")
$ws1.Range("A10").Value = 10
$ws1.Range("B10").Value = 12

# Make "Test Sheet 1" the active sheet/tab and select the newly added cell
$ws1.Activate()
$ws1.Range("B10").Select()

$wb.Save()
